# Chitter Challenge Steps.docx
# Append two new list-paragraphs after "Moving onto signing out." describing
# the sign-out feature work and the decision to DRY up sign in/out into a
# shared module.

$d = $word.ActiveDocument

# Locate the final paragraph of the document ("Moving onto signing out.")
$anchor = $d.Paragraphs.Last
$anchorRange = $anchor.Range
$anchorRange.Collapse(0)            # wdCollapseEnd
$anchorRange.InsertParagraphAfter() # new (still empty) ListParagraph, inherits pPr/rPr

# --- Paragraph 1: two runs sharing identical formatting --------------------
$para1Target = $d.Paragraphs.Last
$para1Range = $para1Target.Range

$para1Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Signing out feature begins with before(:each). I </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">may just go along with this instead of using factory girl and see what happens. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para1Range.InsertXML($para1Xml)

# --- Paragraph 2: single run ------------------------------------------------
$para2Target = $d.Paragraphs.Last
$para2Range = $para2Target.Range

$para2Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Now our sign in method is being used in two separate features and that’s not keeping with DRY principles. I’ll place them in a module and include it in spec helper. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para2Range.InsertXML($para2Xml)

# InsertXML leaves the paragraph it targeted behind (pushed after the newly
# inserted content) as a now-redundant empty paragraph. Remove it by deleting
# the range that spans the preceding paragraph mark through this paragraph's
# end.
$trailing = $d.Paragraphs.Last
$trailingRange = $trailing.Range
$prev = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$cleanupRange = $d.Range($prev.Range.End - 1, $trailingRange.End)
$cleanupRange.Delete()

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
